$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (ALC)
$ws.Range("H2").Value = 164.35
$ws.Range("I2").Value = 172.61111
$ws.Range("K2").Value = 172.61111
$ws.Range("M2").Value = -59.61111

# Row 11 (ALC)
$ws.Range("H11").Value = 100.64286
$ws.Range("I11").Value = 100.64286
$ws.Range("K11").Value = 100.64286
$ws.Range("M11").Value = 39.35714

# Row 40 (ALC)
$ws.Range("H40").Value = 6104.4
$ws.Range("I40").Value = 4008.8
$ws.Range("J40").Value = 8200
$ws.Range("K40").Value = 4008.8
$ws.Range("L40").Value = 8200
$ws.Range("M40").Value = -3833.8
$ws.Range("N40").Value = -8550

# Row 61 (ALC)
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

# Row 64 (ALC)
$ws.Range("H64").Value = 9262.888999999999
$ws.Range("J64").Value = 13400
$ws.Range("L64").Value = 13400
$ws.Range("N64").Value = -13896

# Row 67 (ALC)
$ws.Range("H67").Value = 9262.888999999999
$ws.Range("J67").Value = 13400
$ws.Range("L67").Value = 13400
$ws.Range("N67").Value = -15116

# Row 80 (ALC)
$ws.Range("H80").Value = 393.35135
$ws.Range("I80").Value = 376.25806
$ws.Range("J80").Value = 481.66666
$ws.Range("K80").Value = 1128.77418
$ws.Range("L80").Value = 1444.99998
$ws.Range("M80").Value = -130.7741799999999
$ws.Range("N80").Value = -3440.99998

# Row 83 (ALC)
$ws.Range("H83").Value = 393.35135
$ws.Range("I83").Value = 376.25806
$ws.Range("J83").Value = 481.66666
$ws.Range("K83").Value = 3386.32254
$ws.Range("L83").Value = 4334.99994
$ws.Range("M83").Value = 1605.67746
$ws.Range("N83").Value = -14318.99994

# Row 106 (ALC)
$ws.Range("H106").Value = 4548.6665
$ws.Range("I106").Value = 4358.4
$ws.Range("K106").Value = 4358.4
$ws.Range("M106").Value = -3727.4

# Row 107 (ALC)
$ws.Range("H107").Value = 529.4737
$ws.Range("I107").Value = 529.4737
$ws.Range("K107").Value = 529.4737
$ws.Range("M107").Value = 1390.5263

# Row 116 (ALC)
$ws.Range("H116").Value = 4400
$ws.Range("J116").Value = 5900
$ws.Range("L116").Value = 5900
$ws.Range("N116").Value = -12784

# Row 135 (ALC)
$ws.Range("H135").Value = 1166.6666
$ws.Range("I135").Value = 1519
$ws.Range("J135").Value = 109.666664
$ws.Range("K135").Value = 13671
$ws.Range("L135").Value = 986.9999759999999
$ws.Range("M135").Value = -11136
$ws.Range("N135").Value = -6056.999976

# Row 141 (ALC)
$ws.Range("H141").Value = 7497.231
$ws.Range("I141").Value = 6872.0835
$ws.Range("K141").Value = 20616.2505
$ws.Range("M141").Value = -15436.2505

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 7784.5625
$ws.Range("I32").Value = 7990.5483
$ws.Range("K32").Value = 7990.5483
$ws.Range("M32").Value = -7703.5483

# Row 45 (ARM)
$ws.Range("H45").Value = 3274.386
$ws.Range("J45").Value = 5220.579
$ws.Range("L45").Value = 5220.579
$ws.Range("N45").Value = -5974.579

# Row 63 (ARM)
$ws.Range("H63").Value = 2121.5
$ws.Range("I63").Value = 1900
$ws.Range("J63").Value = 3229
$ws.Range("K63").Value = 1900
$ws.Range("L63").Value = 3229
$ws.Range("M63").Value = -1214
$ws.Range("N63").Value = -4601

# Row 66 (ARM)
$ws.Range("H66").Value = 2121.5
$ws.Range("I66").Value = 1900
$ws.Range("J66").Value = 3229
$ws.Range("K66").Value = 9500
$ws.Range("L66").Value = 16145
$ws.Range("M66").Value = -6068
$ws.Range("N66").Value = -23009

# Row 74 (ARM)
$ws.Range("H74").Value = 6818.375
$ws.Range("I74").Value = 4977.3335
$ws.Range("J74").Value = 9185.429
$ws.Range("K74").Value = 4977.3335
$ws.Range("L74").Value = 9185.429
$ws.Range("M74").Value = -4103.3335
$ws.Range("N74").Value = -10933.429

# Row 77 (ARM)
$ws.Range("H77").Value = 6818.375
$ws.Range("I77").Value = 4977.3335
$ws.Range("J77").Value = 9185.429
$ws.Range("K77").Value = 24886.6675
$ws.Range("L77").Value = 45927.145
$ws.Range("M77").Value = -20518.6675
$ws.Range("N77").Value = -54663.145

# Row 110 (ARM)
$ws.Range("H110").Value = 7949.8
$ws.Range("I110").Value = 5785.4287
$ws.Range("K110").Value = 5785.4287
$ws.Range("M110").Value = -3740.4287

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 2491.75
$ws.Range("I20").Value = 2611.9167
$ws.Range("J20").Value = 2311.5
$ws.Range("K20").Value = 2611.9167
$ws.Range("L20").Value = 2311.5
$ws.Range("M20").Value = -2364.9167
$ws.Range("N20").Value = -2805.5

# Row 86 (BSM)
$ws.Range("H86").Value = 40003684
$ws.Range("I86").Value = 5977.5
$ws.Range("J86").Value = 66668820
$ws.Range("K86").Value = 5977.5
$ws.Range("L86").Value = 66668820
$ws.Range("M86").Value = -4854.5
$ws.Range("N86").Value = -66671066

# Row 89 (BSM)
$ws.Range("H89").Value = 40003684
$ws.Range("I89").Value = 5977.5
$ws.Range("J89").Value = 66668820
$ws.Range("K89").Value = 29887.5
$ws.Range("L89").Value = 333344100
$ws.Range("M89").Value = -24271.5
$ws.Range("N89").Value = -333355332

# Row 107 (BSM)
$ws.Range("H107").Value = 2888.303
$ws.Range("I107").Value = 2389.9092
$ws.Range("J107").Value = 3885.0908
$ws.Range("K107").Value = 2389.9092
$ws.Range("L107").Value = 3885.0908
$ws.Range("M107").Value = -469.9092000000001
$ws.Range("N107").Value = -7725.0908

# Row 130 (BSM)
$ws.Range("H130").Value = 25000
$ws.Range("I130").Value = 25000
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 25000
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -19980
$ws.Range("N130").ClearContents()

# Row 134 (BSM)
$ws.Range("H134").Value = 4970.3335
$ws.Range("I134").Value = 5112.0835
$ws.Range("J134").Value = 3836.3333
$ws.Range("K134").Value = 15336.2505
$ws.Range("L134").Value = 11508.9999
$ws.Range("M134").Value = -12801.2505
$ws.Range("N134").Value = -16578.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 4469.515
$ws.Range("I31").Value = 4145.1333
$ws.Range("J31").Value = 4739.8335
$ws.Range("K31").Value = 4145.1333
$ws.Range("L31").Value = 4739.8335
$ws.Range("M31").Value = -3850.1333
$ws.Range("N31").Value = -5329.8335

# Row 34 (CRP)
$ws.Range("H34").Value = 4469.515
$ws.Range("I34").Value = 4145.1333
$ws.Range("J34").Value = 4739.8335
$ws.Range("K34").Value = 4145.1333
$ws.Range("L34").Value = 4739.8335
$ws.Range("M34").Value = -3943.1333
$ws.Range("N34").Value = -5143.8335

# Row 58 (CRP)
$ws.Range("H58").Value = 7753.769
$ws.Range("I58").Value = 3471.2856
$ws.Range("J58").Value = 12750
$ws.Range("K58").Value = 3471.2856
$ws.Range("L58").Value = 12750
$ws.Range("M58").Value = -3268.2856
$ws.Range("N58").Value = -13156

# Row 62 (CRP)
$ws.Range("H62").Value = 5670
$ws.Range("I62").Value = 4005
$ws.Range("J62").Value = 6502.5
$ws.Range("K62").Value = 4005
$ws.Range("L62").Value = 6502.5
$ws.Range("M62").Value = -3381
$ws.Range("N62").Value = -7750.5

# Row 65 (CRP)
$ws.Range("H65").Value = 5670
$ws.Range("I65").Value = 4005
$ws.Range("J65").Value = 6502.5
$ws.Range("K65").Value = 20025
$ws.Range("L65").Value = 32512.5
$ws.Range("M65").Value = -16905
$ws.Range("N65").Value = -38752.5

# Row 99 (CRP)
$ws.Range("H99").Value = 4607.5
$ws.Range("I99").Value = 3762.375
$ws.Range("J99").Value = 5959.7
$ws.Range("K99").Value = 3762.375
$ws.Range("L99").Value = 5959.7
$ws.Range("M99").Value = -2264.375
$ws.Range("N99").Value = -8955.700000000001

# Row 122 (CRP)
$ws.Range("H122").Value = 3964.6
$ws.Range("J122").Value = 5366.6665
$ws.Range("L122").Value = 16099.9995
$ws.Range("N122").Value = -20999.9995

# Row 126 (CRP)
$ws.Range("H126").Value = 4607.5
$ws.Range("I126").Value = 3762.375
$ws.Range("J126").Value = 5959.7
$ws.Range("K126").Value = 11287.125
$ws.Range("L126").Value = 17879.1
$ws.Range("M126").Value = -8817.125
$ws.Range("N126").Value = -22819.1

# Row 136 (CRP)
$ws.Range("H136").Value = 7753.769
$ws.Range("I136").Value = 3471.2856
$ws.Range("J136").Value = 12750
$ws.Range("K136").Value = 10413.8568
$ws.Range("L136").Value = 38250
$ws.Range("M136").Value = -7863.856800000001
$ws.Range("N136").Value = -43350

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 14142.714
$ws.Range("I70").Value = 4499.5
$ws.Range("K70").Value = 4499.5
$ws.Range("M70").Value = -4229.5

# Row 73 (GSM)
$ws.Range("H73").Value = 14142.714
$ws.Range("I73").Value = 4499.5
$ws.Range("K73").Value = 4499.5
$ws.Range("M73").Value = -3563.5

# Row 107 (GSM)
$ws.Range("H107").Value = 464.4
$ws.Range("J107").Value = 502.33334
$ws.Range("L107").Value = 502.33334
$ws.Range("N107").Value = -4342.33334

# Row 113 (GSM)
$ws.Range("H113").Value = 31224.5
$ws.Range("I113").Value = 58449.5
$ws.Range("J113").Value = 3999.5
$ws.Range("K113").Value = 58449.5
$ws.Range("L113").Value = 3999.5
$ws.Range("M113").Value = -56279.5
$ws.Range("N113").Value = -8339.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2564.8572
$ws.Range("I7").Value = 2492.923
$ws.Range("K7").Value = 2492.923
$ws.Range("M7").Value = -2380.923

# Row 22 (LTW)
$ws.Range("H22").Value = 873.6667
$ws.Range("I22").Value = 945
$ws.Range("J22").Value = 864.75
$ws.Range("K22").Value = 945
$ws.Range("L22").Value = 864.75
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -1454.75

# Row 27 (LTW)
$ws.Range("H27").Value = 873.6667
$ws.Range("I27").Value = 945
$ws.Range("J27").Value = 864.75
$ws.Range("K27").Value = 945
$ws.Range("L27").Value = 864.75
$ws.Range("M27").Value = -838
$ws.Range("N27").Value = -1078.75

# Row 40 (LTW)
$ws.Range("H40").Value = 3791.7856
$ws.Range("I40").Value = 2908.9
$ws.Range("K40").Value = 2908.9
$ws.Range("M40").Value = -2772.9

# Row 68 (LTW)
$ws.Range("H68").Value = 12442.25
$ws.Range("I68").Value = 11973.429
$ws.Range("K68").Value = 11973.429
$ws.Range("M68").Value = -11224.429

# Row 71 (LTW)
$ws.Range("H71").Value = 12442.25
$ws.Range("I71").Value = 11973.429
$ws.Range("K71").Value = 59867.145
$ws.Range("M71").Value = -56123.145

# Row 126 (LTW)
$ws.Range("H126").Value = 2564.8572
$ws.Range("I126").Value = 2492.923
$ws.Range("K126").Value = 7478.768999999999
$ws.Range("M126").Value = -5008.768999999999

# Row 140 (LTW)
$ws.Range("H140").Value = 94275
$ws.Range("J140").Value = 94275
$ws.Range("L140").Value = 94275
$ws.Range("N140").Value = -104635

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (WVR)
$ws.Range("H54").Value = 32000
$ws.Range("J54").Value = 32000
$ws.Range("L54").Value = 32000
$ws.Range("N54").Value = -33040
